$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "Completed till lesson 5.13"
$ws.Range("B11").Value = "D17"
$ws.Range("C10").Copy($ws.Range("C11"))
$ws.Range("C11").Value = 43851

$ws.Range("B11").Select()
